# Update "想去人数" (F column) values on the 展览, 演出 and 全部类型 sheets,
# reflecting the refreshed stats captured in the latest data pull.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 273
$ws1.Range("F5").Value  = 279
$ws1.Range("F7").Value  = 102
$ws1.Range("F8").Value  = 264
$ws1.Range("F11").Value = 38
$ws1.Range("F12").Value = 117
$ws1.Range("F13").Value = 2462
$ws1.Range("F14").Value = 49
$ws1.Range("F15").Value = 24
$ws1.Range("F19").Value = 534
$ws1.Range("F20").Value = 581
$ws1.Range("F21").Value = 175
$ws1.Range("F22").Value = 89
$ws1.Range("F24").Value = 52
$ws1.Range("F25").Value = 2050
$ws1.Range("F26").Value = 4153
$ws1.Range("F28").Value = 66
$ws1.Range("F29").Value = 466
$ws1.Range("F30").Value = 1216
$ws1.Range("F31").Value = 238
$ws1.Range("F32").Value = 2121
$ws1.Range("F35").Value = 66
$ws1.Range("F37").Value = 294
$ws1.Range("F39").Value = 718
$ws1.Range("F40").Value = 8
$ws1.Range("F41").Value = 442
$ws1.Range("F43").Value = 429

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 43

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 273
$ws4.Range("F5").Value  = 279
$ws4.Range("F7").Value  = 102
$ws4.Range("F8").Value  = 264
$ws4.Range("F11").Value = 38
$ws4.Range("F12").Value = 117
$ws4.Range("F13").Value = 2462
$ws4.Range("F14").Value = 49
$ws4.Range("F15").Value = 24
$ws4.Range("F17").Value = 44
$ws4.Range("F20").Value = 534
$ws4.Range("F21").Value = 581
$ws4.Range("F22").Value = 175
$ws4.Range("F23").Value = 89
$ws4.Range("F25").Value = 52
$ws4.Range("F26").Value = 2050
$ws4.Range("F27").Value = 4153
$ws4.Range("F29").Value = 66
$ws4.Range("F30").Value = 466
$ws4.Range("F31").Value = 1216
$ws4.Range("F32").Value = 238
$ws4.Range("F33").Value = 2121
$ws4.Range("F36").Value = 66
$ws4.Range("F38").Value = 294
$ws4.Range("F40").Value = 718
$ws4.Range("F41").Value = 8
$ws4.Range("F42").Value = 442
$ws4.Range("F44").Value = 429
